# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-08 08:24:08
# Applies the updated "Recorded By" orderings, refreshed session/coverage
# statistics, and the newly-recorded PATHOLOGY LAB/MUSEUM C1 session (row 25)
# to the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write literal text into a cell while preserving that cell's current
# style (so things like "79.3%" stay as plain text instead of being
# reinterpreted as a percentage number by Excel's input parser).
# ---------------------------------------------------------------------------
function Set-LiteralText {
    param($Cell, [string]$Text)

    $helper = $ws.Range("ZZ1")
    $helper.Value = "'" + $Text
    $helper.Copy()
    $Cell.PasteSpecial(-4163)  # xlPasteValues
    $helper.Clear()
}

# ---------------------------------------------------------------------------
# Row 2 - ANATOMY session 1 : reorder "Recorded By" list
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = "Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 3 - ANATOMY session 2 : reorder "Recorded By" list
# ---------------------------------------------------------------------------
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 4 - ANATOMY session 3 : reorder "Recorded By" list
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 5 - ANATOMY session 4 : reorder "Recorded By" list
# ---------------------------------------------------------------------------
$ws.Range("G5").Value = "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 6 - Class Statistics : Recorded Sessions 22 -> 23
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 23

# ---------------------------------------------------------------------------
# Row 7 - BIOCHEMISTRY LAB/CBL session 1 : reorder "Recorded By" list,
#         Class Statistics : Missing Sessions 3 -> 2
# ---------------------------------------------------------------------------
$ws.Range("G7").Value = "AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
$ws.Range("L7").Value = 2

# ---------------------------------------------------------------------------
# Row 9 - Class Statistics : Coverage % 75.9% -> 79.3%
# ---------------------------------------------------------------------------
Set-LiteralText ($ws.Range("L9")) "79.3%"

# ---------------------------------------------------------------------------
# Row 10 - Class Statistics : Average Attendance % 27.3% -> 27.1%
# ---------------------------------------------------------------------------
Set-LiteralText ($ws.Range("L10")) "27.1%"

# ---------------------------------------------------------------------------
# Row 12 - MICROBIOLOGY session 1 : reorder "Recorded By" list
# ---------------------------------------------------------------------------
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 15 - Group Statistics summary row (mirrors rows 6/7/9/10 for Year 2 C1)
# ---------------------------------------------------------------------------
$ws.Range("O15").Value = 23
$ws.Range("P15").Value = 2
Set-LiteralText ($ws.Range("R15")) "79.3%"
Set-LiteralText ($ws.Range("S15")) "27.1%"

# ---------------------------------------------------------------------------
# Row 24 - PATHOLOGY LAB/MUSEUM session 2 : reorder "Recorded By" list
# ---------------------------------------------------------------------------
$ws.Range("G24").Value = "Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 25 - PATHOLOGY LAB/MUSEUM session 2 (08/12/2025) : this session has now
# been recorded. Re-style the row from the "Not Recorded" (pink) look to the
# "Recorded" (green) look used elsewhere, fill in the recorder and student
# count, and flip the status.
# ---------------------------------------------------------------------------
$ws.Range("A24:I24").Copy()
$ws.Range("A25:I25").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G25").Value = "menna-allah.gamil@med.asu.edu.eg"
$ws.Range("H25").Value = "60/251"
$ws.Range("I25").Value = "Recorded"

# ---------------------------------------------------------------------------
# Row 27 - PHARMACOLOGY session 2 : reorder "Recorded By" list
# ---------------------------------------------------------------------------
$ws.Range("G27").Value = "nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 28 - PHYSIOLOGY session 1 : reorder "Recorded By" list
# ---------------------------------------------------------------------------
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 30 - PHYSIOLOGY session 3 : reorder "Recorded By" list
# ---------------------------------------------------------------------------
$ws.Range("G30").Value = "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
